# CAN2-31 tx priority reqs generated ready for review cr1
#
# This script edits CAN_Controller_Requirements_ColinFritz.xlsx:
#  - Adds a second "Priority" requirement row (TxPL_PRIORITY_02) to the
#    Tx_Priority_Logic sheet, right after the existing TxPL_PRIORITY_01 row.
#  - Extends the blank filler rows below that table by two rows.
#  - Re-merges the section label columns to account for the new row.
#  - Makes Tx_Priority_Logic the active / selected sheet (it was
#    Acceptance_Filter before).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tx_Priority_Logic")
$ws2 = $wb.Worksheets.Item("Acceptance_Filter")

# ---------------------------------------------------------------------
# 1. Make room for the new requirement row: push the existing blank
#    spacer row (old row 10) and everything below it down by one row.
# ---------------------------------------------------------------------
$ws1.Range("A10:D10").Insert(-4121)

# ---------------------------------------------------------------------
# 2. Un-merge / re-merge the section label column now that the Priority
#    block spans rows 2-10 instead of 2-9.
# ---------------------------------------------------------------------
$ws1.Range("A2:A9").UnMerge()
$ws1.Range("A2:A10").Merge()

# ---------------------------------------------------------------------
# 3. Row 9 is no longer the last row of the Priority block, so it picks
#    up the same "interior row" borders used by every other interior
#    row of a multi-row block (copy from the Send Data block, which has
#    the identical shape).
# ---------------------------------------------------------------------

# A9 becomes an interior-row cell (same look as A3/A4/A6/A7/A8).
$ws1.Range("A3").Copy()
$ws1.Range("A9").PasteSpecial(-4122)

# B9 becomes the top cell of the new B9:B10 merge.
$ws1.Range("A2").Copy()
$ws1.Range("B9").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Row 10 becomes the new bottom row of the Priority block: same
#    border treatment the old row 9 used to have (bottom border).
# ---------------------------------------------------------------------
$ws1.Range("A9").Copy()
$ws1.Range("A10").PasteSpecial(-4122)
$ws1.Range("A9").Copy()
$ws1.Range("B10").PasteSpecial(-4122)

$ws1.Range("C1").Copy()
$ws1.Range("C10").PasteSpecial(-4122)
$ws1.Range("C1").Copy()
$ws1.Range("D10").PasteSpecial(-4122)

$ws1.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 5. Merge the Priority label cell across the two requirement rows and
#    fill in the new requirement's text.
# ---------------------------------------------------------------------
$ws1.Range("B9:B10").Merge()

$ws1.Range("C10").Value = "TxPL_PRIORITY_02"
$ws1.Range("D10").Value = "The module shall implement full and empty flags for indicating the occupancy of the internal register.  "

# ---------------------------------------------------------------------
# 6. Restore the blank spacer rows below the table (rows 11-14 keep the
#    old vertically-centered style) and append two more blank rows so
#    the sheet ends at row 17, matching the extended table.
# ---------------------------------------------------------------------
$ws1.Range("A15:A16").Insert(-4121)
$ws1.Range("A17").Value = ""

for ($r = 11; $r -le 17; $r++) {
    $ws1.Range("A11").Copy()
    $ws1.Range("A$r").PasteSpecial(-4122)
}
$ws1.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 7. Tx_Priority_Logic is now the sheet the user is looking at.
# ---------------------------------------------------------------------
$ws1.Select()
$ws1.Range("D18").Select()
